{"js": "// The document's last six \"Listaszerbekezds\" items track progress on the\n// \"5 allek\u00e9rdez\u00e9st haszn\u00e1l\u00f3 lek\u00e9rdez\u00e9s\" (5 subqueries) exercise. This edit\n// replaces the leftover \"TODO\" placeholder paragraph with real progress,\n// rotates the per-item notes forward, drops the now-finished\n// \"Lak\u00f3_Legfiatalabb_...(ALL)\" line, and moves the \"_GoBack\" bookmark\n// (Word's \"last edited location\" marker) to the paragraph that now holds\n// the newest edit.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Locate the anchor paragraphs by their current text so this keeps working\n// even if unrelated paragraphs shift around.\nparagraphs.items.forEach((p) => p.load(\"text\"));\nawait context.sync();\n\nconst todoPara = paragraphs.items.find((p) => p.text.startsWith(\"TODO\"));\nconst minPara = paragraphs.items.find((p) => p.text.indexOf(\"Lak\u00f3_legfiatalabb_\") === 0 && p.text.indexOf(\"minnel\") !== -1);\nconst maxPara = paragraphs.items.find((p) => p.text.indexOf(\"Lak\u00f3_legid\u0151sebb_\") === 0 && p.text.indexOf(\"maxxal\") !== -1);\nconst allPara = paragraphs.items.find((p) => p.text.indexOf(\"Lak\u00f3_Legfiatalabb_\") === 0 && p.text.indexOf(\"(ALL)\") !== -1);\nconst anyPara = paragraphs.items.find((p) => p.text.indexOf(\"Lak\u00f3_legid\u0151sebb_\") === 0 && p.text.indexOf(\"(ANY)\") !== -1);\n\nconst W_DOC_OPEN = '<?xml version=\"1.0\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>';\nconst W_DOC_CLOSE = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>';\n\nfunction listParaOoxml(innerRunsXml, opts) {\n  opts = opts || {};\n  const bookmark = opts.bookmark\n    ? '<w:bookmarkStart w:id=\"0\" w:name=\"' + opts.bookmark + '\"/><w:bookmarkEnd w:id=\"0\"/>'\n    : \"\";\n  return W_DOC_OPEN +\n    '<w:p><w:pPr><w:pStyle w:val=\"Listaszerbekezds\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"6\"/></w:numPr></w:pPr>' +\n    bookmark + innerRunsXml + '</w:p>' +\n    W_DOC_CLOSE;\n}\n\n// 0) Remove the old \"_GoBack\" bookmark first (it currently lives on the\n//    \"(ANY)\" paragraph) so re-inserting it below doesn't leave two bookmarks\n//    with the same name alive at once.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// 1) Old \"TODO: Form\u00e1zott ki\u00edr\u00e1s\" paragraph -> \"Lak\u00f3_legfiatalabb_maxxal\",\n//    now carrying the numbered-list formatting and the _GoBack bookmark.\nconst todoRuns =\n  '<w:r><w:t>Lak\u00f3_legfiatalabb_</w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/><w:r><w:t>maxxal</w:t></w:r><w:proofErr w:type=\"spellEnd\"/>';\ntodoPara.getRange().insertOoxml(listParaOoxml(todoRuns, { bookmark: \"_GoBack\" }), Word.InsertLocation.replace);\n\n// 2) \"Lak\u00f3_legfiatalabb_minnel\" -> \"Lak\u00f3_legid\u0151sebb_minnel \"\nconst minRuns =\n  '<w:r><w:t>Lak\u00f3_legid\u0151sebb_</w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/><w:r><w:t>minnel</w:t></w:r><w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>';\nminPara.getRange().insertOoxml(listParaOoxml(minRuns), Word.InsertLocation.replace);\n\n// 3) \"Lak\u00f3_legid\u0151sebb_maxxal\" -> \"Lak\u00f3_Legfiatalabb_All (ALL)\"\nconst maxRuns =\n  '<w:r><w:t>L</w:t></w:r>' +\n  '<w:r><w:t>ak\u00f3_Legfiatalabb_</w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/><w:r><w:t>All</w:t></w:r><w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r><w:t xml:space=\"preserve\"> (</w:t></w:r>' +\n  '<w:r><w:t>ALL</w:t></w:r>' +\n  '<w:r><w:t>)</w:t></w:r>';\nmaxPara.getRange().insertOoxml(listParaOoxml(maxRuns), Word.InsertLocation.replace);\n\n// 4) Drop the finished \"Lak\u00f3_Legfiatalabb_Allekerdezessel (ALL)\" paragraph.\n//    (The \"(ANY)\" paragraph's text is unchanged; its _GoBack bookmark was\n//    already relocated to the first paragraph in step 0/1 above.)\nallPara.delete();\n\nawait context.sync();\n", "ps1": "# The document's last six \"Listaszerbekezds\" items track progress on the\n# \"5 allek\u00e9rdez\u00e9st haszn\u00e1l\u00f3 lek\u00e9rdez\u00e9s\" (5 subqueries) exercise. This edit\n# replaces the leftover \"TODO\" placeholder paragraph with real progress,\n# rotates the per-item notes forward, drops the now-finished\n# \"Lak\u00f3_Legfiatalabb_...(ALL)\" line, and moves the \"_GoBack\" bookmark\n# (Word's \"last edited location\" marker) to the paragraph that now holds\n# the newest edit.\n\n$d = $word.ActiveDocument\n\nfunction Find-ParaByText($doc, $prefix, $contains) {\n    foreach ($p in $doc.Paragraphs) {\n        $t = $p.Range.Text\n        if ($t.StartsWith($prefix) -and $t.Contains($contains)) {\n            return $p\n        }\n    }\n    return $null\n}\n\n$todoPara = Find-ParaByText $d \"TODO\" \"\"\n$minPara  = Find-ParaByText $d \"Lak\u00f3_legfiatalabb_\" \"minnel\"\n$maxPara  = Find-ParaByText $d \"Lak\u00f3_legid\u0151sebb_\" \"maxxal\"\n$allPara  = Find-ParaByText $d \"Lak\u00f3_Legfiatalabb_\" \"(ALL)\"\n\n$xmlOpen = '<?xml version=\"1.0\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>'\n$xmlClose = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n$numPr = '<w:pPr><w:pStyle w:val=\"Listaszerbekezds\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"6\"/></w:numPr></w:pPr>'\n$bookmark = '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/>'\n\n# 0) Remove the old \"_GoBack\" bookmark first (it currently lives on the\n#    \"(ANY)\" paragraph) so re-inserting it below doesn't leave two bookmarks\n#    with the same name alive at once.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n# 1) Old \"TODO: Form\u00e1zott ki\u00edr\u00e1s\" paragraph -> \"Lak\u00f3_legfiatalabb_maxxal\",\n#    now carrying the numbered-list formatting and the _GoBack bookmark.\n$todoRuns = '<w:r><w:t>Lak\u00f3_legfiatalabb_</w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>maxxal</w:t></w:r><w:proofErr w:type=\"spellEnd\"/>'\n$todoXml = $xmlOpen + '<w:p>' + $numPr + $bookmark + $todoRuns + '</w:p>' + $xmlClose\n$null = $todoPara.Range.InsertXML($todoXml)\n\n# 2) \"Lak\u00f3_legfiatalabb_minnel\" -> \"Lak\u00f3_legid\u0151sebb_minnel \"\n$minRuns = '<w:r><w:t>Lak\u00f3_legid\u0151sebb_</w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>minnel</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> </w:t></w:r>'\n$minXml = $xmlOpen + '<w:p>' + $numPr + $minRuns + '</w:p>' + $xmlClose\n$null = $minPara.Range.InsertXML($minXml)\n\n# 3) \"Lak\u00f3_legid\u0151sebb_maxxal\" -> \"Lak\u00f3_Legfiatalabb_All (ALL)\"\n$maxRuns = '<w:r><w:t>L</w:t></w:r><w:r><w:t>ak\u00f3_Legfiatalabb_</w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>All</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> (</w:t></w:r><w:r><w:t>ALL</w:t></w:r><w:r><w:t>)</w:t></w:r>'\n$maxXml = $xmlOpen + '<w:p>' + $numPr + $maxRuns + '</w:p>' + $xmlClose\n$null = $maxPara.Range.InsertXML($maxXml)\n\n# 4) Drop the finished \"Lak\u00f3_Legfiatalabb_Allekerdezessel (ALL)\" paragraph.\n#    (The \"(ANY)\" paragraph's text is unchanged; its _GoBack bookmark was\n#    already relocated to the first paragraph in step 0/1 above.)\n$allPara.Range.Delete()\n"}
